$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.971.84'
$ws.Range("E2").Value = '  -3.65%  '
$ws.Range("D3").Value = '3.346.61'
$ws.Range("E3").Value = '  -4.59%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '''183.04'
$ws.Range("E5").Value = '  -8.98%  '
$ws.Range("D6").Value = '''534.49'
$ws.Range("E6").Value = '  -3.32%  '
$ws.Range("D7").Value = '''0.608'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = '3.344.07'
$ws.Range("E8").Value = '  -4.47%  '
$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '''0.621'
$ws.Range("E10").Value = '  -5.18%  '
$ws.Range("D11").Value = '''59.51'
$ws.Range("E11").Value = '  -6.60%  '
$ws.Range("E12").Value = '  -4.90%  '
$ws.Range("E13").Value = '  -2.48%  '
$ws.Range("D14").Value = '''9.24'
$ws.Range("E14").Value = '  -6.57%  '
$ws.Range("D15").Value = '3.867.86'
$ws.Range("E15").Value = '  -4.38%  '
$ws.Range("D16").Value = '3.342.81'
$ws.Range("E16").Value = '  -4.23%  '
$ws.Range("E17").Value = '  -4.18%  '
$ws.Range("D18").Value = '''17.78'
$ws.Range("E18").Value = '  -3.46%  '
$ws.Range("D19").Value = '65.012.49'
$ws.Range("E19").Value = '  -3.10%  '
$ws.Range("D20").Value = '''11.30'
$ws.Range("E20").Value = '  -4.23%  '
$ws.Range("E21").Value = '  -4.96%  '
$ws.Range("D22").Value = '''379.18'
$ws.Range("E22").Value = '  -2.96%  '
$ws.Range("D23").Value = '''3.86'
$ws.Range("E23").Value = '  -3.65%  '
$ws.Range("D24").Value = '''11.43'
$ws.Range("E24").Value = '  -6.84%  '
$ws.Range("D25").Value = '''81.44'
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("E26").Value = '  +3.20%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '''6.11'
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = '''2.72'
$ws.Range("E28").Value = '  -3.05%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''11.64'
$ws.Range("E29").Value = '  -5.01%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").Value = '''8.51'
$ws.Range("E30").Value = '  -3.59%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '''29.27'
$ws.Range("E31").Value = '  -5.60%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '''660.11'
$ws.Range("E32").Value = '  -2.56%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '''6.81'
$ws.Range("E33").Value = '  -2.85%  '
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").Value = '''11.41'
$ws.Range("E34").Value = '  -3.07%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.108'
$ws.Range("E35").Value = '  -2.33%  '
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").Value = '''59.81'
$ws.Range("E36").Value = '  -6.48%  '
$ws.Range("D37").Value = '''0.399'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").Value = '''37.32'
$ws.Range("E39").Value = '  -3.93%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0717'
$ws.Range("E40").Value = '  +6.31%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''0.996'
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.127'
$ws.Range("E42").Value = '  -2.48%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.946.91'
$ws.Range("E43").Value = '  -4.27%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '''2.56'
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").Value = '''2.75'
$ws.Range("E45").Value = '  -7.52%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0404'
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '''2.67'
$ws.Range("E47").Value = '  -3.66%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '''3.12'
$ws.Range("E48").Value = '  +8.20%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '''2.84'
$ws.Range("E49").Value = '  +8.71%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '''0.127'
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '''2.57'
$ws.Range("E51").Value = '  -4.92%  '
